# Update cryptos list (Price + Volume(1h) columns), as produced by the
# scheduled "Updated cryptos list" GitHub Actions run.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row -> new Price (column D, $null means unchanged) / new Volume(1h) (column E)
$updates = @(
    @{ Row = 2;  D = "26.719.85";  E = "  +1.83%  " },
    @{ Row = 3;  D = "1.623.95";   E = "  +2.23%  " },
    @{ Row = 4;  D = $null;        E = "  -0.01%  " },
    @{ Row = 5;  D = "214.62";     E = "  +1.21%  " },
    @{ Row = 6;  D = $null;        E = "  +0.86%  " },
    @{ Row = 7;  D = $null;        E = "  +0.02%  " },
    @{ Row = 8;  D = $null;        E = "  +0.38%  " },
    @{ Row = 9;  D = $null;        E = "  +0.64%  " },
    @{ Row = 10; D = "19.41";      E = "  +0.47%  " },
    @{ Row = 11; D = $null;        E = "  +1.24%  " },
    @{ Row = 12; D = "1.853.11";   E = "  +2.28%  " },
    @{ Row = 13; D = "1.617.81";   E = "  +1.84%  " },
    @{ Row = 14; D = "4.06";       E = "  +1.28%  " },
    @{ Row = 15; D = "65.05";      E = "  +1.26%  " },
    @{ Row = 16; D = "0.513";      E = "  -1.23%  " },
    @{ Row = 17; D = "26.742.92";  E = "  +1.92%  " },
    @{ Row = 18; D = "234.89";     E = "  +10.27%  " },
    @{ Row = 19; D = "7.73";       E = "  +4.82%  " },
    @{ Row = 20; D = $null;        E = "  +0.35%  " },
    @{ Row = 21; D = "1.00";       E = "  +0.00%  " },
    @{ Row = 22; D = "4.40";       E = "  +3.21%  " },
    @{ Row = 23; D = "2.28";       E = "  +5.00%  " },
    @{ Row = 24; D = $null;        E = "  +1.33%  " },
    @{ Row = 25; D = "145.94";     E = "  +1.61%  " },
    @{ Row = 26; D = $null;        E = "  +0.07%  " },
    @{ Row = 27; D = "7.07";       E = "  +0.89%  " },
    @{ Row = 28; D = $null;        E = "  +2.40%  " },
    @{ Row = 29; D = "15.69";      E = "  +3.24%  " },
    @{ Row = 30; D = "0.0498";     E = "  +0.40%  " },
    @{ Row = 31; D = $null;        E = "  +1.27%  " },
    @{ Row = 32; D = $null;        E = "  +1.96%  " },
    @{ Row = 33; D = "1.471.25";   E = "  +9.92%  " },
    @{ Row = 34; D = $null;        E = "  +2.35%  " },
    @{ Row = 35; D = "2.43";       E = "  -0.54%  " },
    @{ Row = 36; D = $null;        E = "  +2.03%  " },
    @{ Row = 37; D = "0.571";      E = "  -1.66%  " },
    @{ Row = 38; D = "0.0167";     E = "  +0.49%  " },
    @{ Row = 39; D = "0.840";      E = "  +1.95%  " },
    @{ Row = 40; D = "5.97";       E = "  +3.84%  " },
    @{ Row = 41; D = $null;        E = "  +0.09%  " },
    @{ Row = 42; D = $null;        E = "  +3.07%  " },
    @{ Row = 43; D = "0.955";      E = "  -0.01%  " },
    @{ Row = 44; D = "1.763.49";   E = "  +2.35%  " },
    @{ Row = 45; D = "0.767";      E = "  -0.04%  " },
    @{ Row = 46; D = "62.06";      E = "  +1.30%  " },
    @{ Row = 47; D = $null;        E = "  +3.20%  " },
    @{ Row = 48; D = $null;        E = "  +2.11%  " },
    @{ Row = 49; D = $null;        E = "  +0.78%  " },
    @{ Row = 50; D = "0.0966";     E = "  -1.10%  " },
    @{ Row = 51; D = "7.51";       E = "  +2.26%  " }
)

foreach ($u in $updates) {
    if ($null -ne $u.D) {
        $priceCell = $ws.Cells.Item($u.Row, 4)
        $text = $u.D
        # The "Price" column stores plain text even for values that look
        # like plain numbers (e.g. "214.62"). Writing such a string via
        # .Value would make Excel auto-convert it to a numeric cell, so
        # force a text entry (leading apostrophe) whenever the new value
        # would otherwise be interpreted as a bare number.
        if ($text -match '^[+-]?\d+(\.\d+)?$') {
            $priceCell.Value = "'" + $text
        } else {
            $priceCell.Value = $text
        }
    }
    $ws.Cells.Item($u.Row, 5).Value = $u.E
}
